$d = $word.ActiveDocument

# Paragraph 11 (Heading1)
$p11 = $d.Paragraphs(11)
$r11 = $p11.Range
$xml11 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t xml:space="preserve">For more information about </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ReleFly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ReleCloud</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, please visit www.relecloud.com or contact press@relecloud.com.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r11.InsertXML($xml11)

# Paragraph 10 (ListParagraph)
$p10 = $d.Paragraphs(10)
$r10 = $p10.Range
$xml10 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>ReleCloud</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is known for its customer-centric culture and its focus on innovation and social responsibility. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ReleCloud</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> has invested heavily in renewable energy, education, health, and philanthropy.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r10.InsertXML($xml10)

# Paragraph 9 (ListParagraph)
$p9 = $d.Paragraphs(9)
$r9 = $p9.Range
$xml9 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>ReleCloud</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> has over 200 million customers worldwide and employs over 500,000 people. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ReleCloud''s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> revenue in 2023 was $500 billion, making it one of the most valuable companies in the world.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r9.InsertXML($xml9)

# Paragraph 8 (ListParagraph)
$p8 = $d.Paragraphs(8)
$r8 = $p8.Range
$xml8 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>ReleCloud</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is a global cloud computing and e-commerce company that offers a wide range of products and services, including online retail, cloud computing, digital streaming, artificial intelligence, and more.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r8.InsertXML($xml8)

# Paragraph 7 (Heading1)
$p7 = $d.Paragraphs(7)
$r7 = $p7.Range
$xml7 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:lastRenderedPageBreak/><w:t>ReleCloud''s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> CEO, Jane Smith, said that </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ReleFly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is a game-changer for the e-commerce industry and the customers. "We are excited to introduce </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ReleFly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, our new delivery drone, that will offer our customers a faster, safer, and greener way to receive their orders. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ReleFly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is the result of years of research and development and reflects our commitment to innovation and customer satisfaction. With </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ReleFly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, we aim to make online shopping more convenient and enjoyable for everyone."</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r7.InsertXML($xml7)

# Paragraph 6 (ListParagraph)
$p6 = $d.Paragraphs(6)
$r6 = $p6.Range
$xml6 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>ReleFly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is designed to reduce carbon emissions and noise pollution by using renewable energy and silent propellers. It can also land and take off vertically, minimizing the need for infrastructure and space.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r6.InsertXML($xml6)

# Paragraph 5 (ListParagraph)
$p5 = $d.Paragraphs(5)
$r5 = $p5.Range
$xml5 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>ReleFly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> uses advanced artificial intelligence and computer vision to navigate complex urban environments and avoid obstacles. It can also communicate with other drones and traffic systems to ensure safe and efficient delivery.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r5.InsertXML($xml5)

# Paragraph 4 (ListParagraph)
$p4 = $d.Paragraphs(4)
$r4 = $p4.Range
$xml4 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>ReleFly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is a lightweight, autonomous, and electric drone that can carry packages up to 5 kg and fly up to 25 km in a single charge.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r4.InsertXML($xml4)

# Paragraph 3 (Heading1)
$p3 = $d.Paragraphs(3)
$r3 = $p3.Range
$xml3 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>ReleCloud</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, a leading cloud computing and e-commerce company, has announced the launch of a new delivery drone that will revolutionize the online shopping experience. The drone, named </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ReleFly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, will be available for customers in summer 2024 and will offer faster, safer, and greener delivery options.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r3.InsertXML($xml3)

# Paragraph 2 (Subtitle)
$p2 = $d.Paragraphs(2)
$r2 = $p2.Range
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Subtitle"/></w:pPr><w:r><w:t>The innovative drone will offer faster, safer, and greener delivery options for customers</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r2.InsertXML($xml2)

# Paragraph 1 (Title)
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>ReleCloud</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Announces the Launch of a New Delivery Drone</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r1.InsertXML($xml1)
